$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.903.64'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.814.99'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.13'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4654'
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3657'
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07361'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8689'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.24'
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = '1.806.95'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.382'
$ws.Range("E13").Value = '  +0.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07104'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.504'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008682'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("E20").Value = '  -0.67%  '
$ws.Range("D21").Value = '26.928.35'
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.55'
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '2.050.49'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.94'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.40'
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.132'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.259'
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.85'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08905'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7565'
$ws.Range("E32").Value = '  +0.53%  '
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.477'
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.911'
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.089'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05275'
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01943'
$ws.Range("E39").Value = '  -1.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.962'
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.178'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5274'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.352'
$ws.Range("E43").Value = '  -2.70%  '
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.428'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("E46").Value = '  -2.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.39'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.15'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.658'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06288'
$ws.Range("E51").Value = '  -0.12%  '
